# Repull data, push all data, mean calculation
# Update column F (dSF) values for several rows to reflect corrected/repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F13").Value = 7
$ws.Range("F15").Value = -2
$ws.Range("F17").Value = 6
$ws.Range("F19").Value = -1
$ws.Range("F24").Value = -5
$ws.Range("F33").Value = -2
$ws.Range("F37").Value = 1
$ws.Range("F38").Value = -5
